$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 320.27777
$ws.Range("I33").Value = 333.92856
$ws.Range("J33").Value = 272.5
$ws.Range("K33").Value = 333.92856
$ws.Range("L33").Value = 272.5
$ws.Range("M33").Value = -104.92856
$ws.Range("N33").Value = -730.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 14999
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 14999
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4188.7856
$ws.Range("I74").Value = 2815.75
$ws.Range("K74").Value = 2815.75
$ws.Range("M74").Value = -1879.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 99594.60000000001
$ws.Range("J75").Value = 99594.60000000001
$ws.Range("L75").Value = 99594.60000000001
$ws.Range("N75").Value = -101466.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4188.7856
$ws.Range("I77").Value = 2815.75
$ws.Range("K77").Value = 14078.75
$ws.Range("M77").Value = -9398.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 99594.60000000001
$ws.Range("J78").Value = 99594.60000000001
$ws.Range("L78").Value = 298783.8
$ws.Range("N78").Value = -308143.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2322.718
$ws.Range("I80").Value = 2283.7144
$ws.Range("J80").Value = 2368.2222
$ws.Range("K80").Value = 6851.1432
$ws.Range("L80").Value = 7104.6666
$ws.Range("M80").Value = -5853.1432
$ws.Range("N80").Value = -9100.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2322.718
$ws.Range("I83").Value = 2283.7144
$ws.Range("J83").Value = 2368.2222
$ws.Range("K83").Value = 20553.4296
$ws.Range("L83").Value = 21313.9998
$ws.Range("M83").Value = -15561.4296
$ws.Range("N83").Value = -31297.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1424.7894
$ws.Range("I100").Value = 516.8182
$ws.Range("K100").Value = 516.8182
$ws.Range("M100").Value = 24.18179999999995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 848
$ws.Range("I103").Value = 844.5
$ws.Range("K103").Value = 2533.5
$ws.Range("M103").Value = -1947.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 143
$ws.Range("I3").Value = 143
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 143
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 255211
$ws.Range("J13").Value = 255211
$ws.Range("L13").Value = 255211
$ws.Range("N13").Value = -255499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1230.4762
$ws.Range("I97").Value = 315.125
$ws.Range("J97").Value = 4159.6
$ws.Range("K97").Value = 315.125
$ws.Range("L97").Value = 4159.6
$ws.Range("M97").Value = 180.875
$ws.Range("N97").Value = -5151.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3459.348
$ws.Range("I122").Value = 1844.3334
$ws.Range("J122").Value = 4497.5713
$ws.Range("K122").Value = 5533.0002
$ws.Range("L122").Value = 13492.7139
$ws.Range("M122").Value = -3083.0002
$ws.Range("N122").Value = -18392.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1743.0294
$ws.Range("I94").Value = 1540.52
$ws.Range("J94").Value = 2305.5557
$ws.Range("K94").Value = 1540.52
$ws.Range("L94").Value = 2305.5557
$ws.Range("M94").Value = -1089.52
$ws.Range("N94").Value = -3207.5557

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3130.7317
$ws.Range("I107").Value = 2036.2858
$ws.Range("K107").Value = 2036.2858
$ws.Range("M107").Value = -116.2858000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 681.875
$ws.Range("I10").Value = 326.66666
$ws.Range("J10").Value = 1747.5
$ws.Range("K10").Value = 326.66666
$ws.Range("L10").Value = 1747.5
$ws.Range("M10").Value = -187.66666
$ws.Range("N10").Value = -2025.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 9606.076999999999
$ws.Range("I105").Value = 1392.4286
$ws.Range("K105").Value = 1392.4286
$ws.Range("M105").Value = 354.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 895.5263
$ws.Range("J107").Value = 1272.2858
$ws.Range("L107").Value = 1272.2858
$ws.Range("N107").Value = -5112.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1630.15
$ws.Range("I134").Value = 1575.3125
$ws.Range("J134").Value = 1849.5
$ws.Range("K134").Value = 4725.9375
$ws.Range("L134").Value = 5548.5
$ws.Range("M134").Value = -2190.9375
$ws.Range("N134").Value = -10618.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2781024.5
$ws.Range("J132").Value = 3925142
$ws.Range("L132").Value = 35326278
$ws.Range("N132").Value = -35331338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 999.6667
$ws.Range("I138").Value = 999.6667
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 2999.0001
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = 2140.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 8044712
$ws.Range("I3").Value = 14747547
$ws.Range("J3").Value = 1310
$ws.Range("K3").Value = 14747547
$ws.Range("L3").Value = 1310
$ws.Range("M3").Value = -14747431
$ws.Range("N3").Value = -1542

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 2445555.2
$ws.Range("I14").Value = 3143785.5
$ws.Range("K14").Value = 3143785.5
$ws.Range("M14").Value = -3143617.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 37810
$ws.Range("I93").Value = 34500
$ws.Range("J93").Value = 41120
$ws.Range("K93").Value = 34500
$ws.Range("L93").Value = 41120
$ws.Range("M93").Value = -32628
$ws.Range("N93").Value = -44864

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3540.5264
$ws.Range("I7").Value = 2962.25
$ws.Range("J7").Value = 4183.0557
$ws.Range("K7").Value = 2962.25
$ws.Range("L7").Value = 4183.0557
$ws.Range("M7").Value = -2850.25
$ws.Range("N7").Value = -4407.0557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3540.5264
$ws.Range("I126").Value = 2962.25
$ws.Range("J126").Value = 4183.0557
$ws.Range("K126").Value = 8886.75
$ws.Range("L126").Value = 12549.1671
$ws.Range("M126").Value = -6416.75
$ws.Range("N126").Value = -17489.1671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 15587
$ws.Range("J51").Value = 20000
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -21020

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 47499.5
$ws.Range("J54").Value = 49999
$ws.Range("L54").Value = 49999
$ws.Range("N54").Value = -51039

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7720.5386
$ws.Range("J62").Value = 7794.727
$ws.Range("L62").Value = 7794.727
$ws.Range("N62").Value = -9042.726999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 7720.5386
$ws.Range("J65").Value = 7794.727
$ws.Range("L65").Value = 38973.635
$ws.Range("N65").Value = -45213.635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 731.5
$ws.Range("I81").Value = 534.8570999999999
$ws.Range("J81").Value = 1190.3334
$ws.Range("K81").Value = 1069.7142
$ws.Range("L81").Value = 2380.6668
$ws.Range("M81").Value = -8.714199999999892
$ws.Range("N81").Value = -4502.6668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 731.5
$ws.Range("I84").Value = 534.8570999999999
$ws.Range("J84").Value = 1190.3334
$ws.Range("K84").Value = 5348.571
$ws.Range("L84").Value = 11903.334
$ws.Range("M84").Value = -44.57099999999991
$ws.Range("N84").Value = -22511.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 33367970
$ws.Range("J122").Value = 1521.5
$ws.Range("L122").Value = 4564.5
$ws.Range("N122").Value = -9464.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 55000
$ws.Range("J125").Value = 55000
$ws.Range("L125").Value = 55000
$ws.Range("N125").Value = -64840
